$d = $word.ActiveDocument

# 1. Update the date
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false,
                         $true, 1, $false, "September 21, 2025", 2)

# 2. Split the address line into city/zip on its own paragraph
$addrRange = $d.Content
$addrRange.Find.Execute("3550 Alden Way #1, San Jose CA 95117", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 0)
$addrRange.Text = "3550 Alden Way #1"
$addrRange.InsertParagraphAfter()
$newParaStart = $addrRange.End + 1
$newPara = $d.Range($newParaStart, $newParaStart).Paragraphs(1)
$newPara.Range.Text = "San Jose, CA 95117"

# 3. Remove the empty "No Spacing" paragraph that immediately follows
#    the "Board of Directors" paragraph.
$bodRange = $d.Content
$bodRange.Find.Execute("Board of Directors", $true, $false, $false, $false, $false,
                        $true, 1, $false, "", 0)
$bodPara = $bodRange.Paragraphs(1)
$afterBodPos = $bodPara.Range.End
$nextPara = $d.Range($afterBodPos, $afterBodPos).Paragraphs(1)
if ($nextPara.Range.Text.TrimEnd("`r", "`a") -eq "" -and $nextPara.Style.NameLocal -eq "No Spacing") {
    $nextPara.Range.Delete()
}
